$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibition) - update "想去人数" (want-to-go count) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 843
$ws1.Range("F4").Value = 1570
$ws1.Range("F5").Value = 722
$ws1.Range("F6").Value = 27

# Sheet "全部类型" (all types) - update "想去人数" (want-to-go count) values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 843
$ws4.Range("F4").Value = 1570
$ws4.Range("F6").Value = 722
$ws4.Range("F7").Value = 27
